# Refresh the COVID-19 "paises" sheet with the newer snapshot of data
# (countries & provincias Spain update). The sheet is kept sorted by
# "Casos totales" (column B) descending, so a handful of rows swap
# country names when their totals leapfrog their neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Junio de 2020 a las 02:10"

# Estados Unidos
$ws.Range("B4").Value = 1859424
$ws.Range("C4").Value = 22254
$ws.Range("D4").Value = 615301
$ws.Range("E4").Value = 1137202
$ws.Range("G4").Value = 726
$ws.Range("H4").Value = 106921

# Brasil
$ws.Range("B5").Value = 529018
$ws.Range("C5").Value = 14169
$ws.Range("D5").Value = 211080
$ws.Range("E5").Value = 287892
$ws.Range("G5").Value = 732
$ws.Range("H5").Value = 30046

# Canada
$ws.Range("B17").Value = 91705
$ws.Range("C17").Value = 758
$ws.Range("D17").Value = 49726
$ws.Range("E17").Value = 34653
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 7326

# Argentina overtakes Israel -> swap rows 44/45
$ws.Range("A44").Value = "Argentina"
$ws.Range("B44").Value = 17415
$ws.Range("C44").Value = 564
$ws.Range("D44").Value = 5521
$ws.Range("E44").Value = 11338
$ws.Range("G44").Value = 17
$ws.Range("H44").Value = 556

$ws.Range("A45").Value = "Israel"
$ws.Range("B45").Value = 17169
$ws.Range("C45").Value = 98
$ws.Range("D45").Value = 14878
$ws.Range("E45").Value = 2006
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 285

# Japon (stays in place, values updated)
$ws.Range("B46").Value = 16884
$ws.Range("C46").Value = 33
$ws.Range("D46").Value = 14502
$ws.Range("E46").Value = 1490
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 892

# Uruguay
$ws.Range("B126").Value = 825
$ws.Range("C126").Value = 2
$ws.Range("D126").Value = 689
$ws.Range("E126").Value = 113
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 23

# Mozambique
$ws.Range("D156").Value = 97
$ws.Range("E156").Value = 155

# Zimbabue jumps ahead of Martinica/Islas Feroe/Mongolia -> rows 159-162 shift
$ws.Range("A159").Value = "Zimbabue"
$ws.Range("B159").Value = 203
$ws.Range("C159").Value = 25
$ws.Range("D159").Value = 29
$ws.Range("E159").Value = 170
$ws.Range("H159").Value = 4

$ws.Range("A160").Value = "Martinica"
$ws.Range("B160").Value = 200
$ws.Range("D160").Value = 98
$ws.Range("E160").Value = 88
$ws.Range("H160").Value = 14

$ws.Range("A161").Value = "Islas Feroe"
$ws.Range("B161").Value = 187
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 187
$ws.Range("E161").Value = 0

$ws.Range("A162").Value = "Mongolia"
$ws.Range("B162").Value = 185
$ws.Range("C162").Value = 6
$ws.Range("D162").Value = 44
$ws.Range("E162").Value = 141
$ws.Range("H162").Value = 0

# Gibraltar (stays) then Libia overtakes Guadalupe -> rows 164/165 swap
$ws.Range("A164").Value = "Libia"
$ws.Range("B164").Value = 168
$ws.Range("C164").Value = 12
$ws.Range("D164").Value = 52
$ws.Range("E164").Value = 111
$ws.Range("H164").Value = 5

$ws.Range("A165").Value = "Guadalupe"
$ws.Range("B165").Value = 162
$ws.Range("D165").Value = 138
$ws.Range("E165").Value = 10
$ws.Range("H165").Value = 14

# Guyana / Islas Caimanes (stay, values updated)
$ws.Range("B167").Value = 150
$ws.Range("C167").Value = 9
$ws.Range("D167").Value = 75
$ws.Range("E167").Value = 74

# Bermudas overtakes Brunei -> rows 168/169 swap
$ws.Range("A168").Value = "Bermudas"
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 112
$ws.Range("E168").Value = 20
$ws.Range("H168").Value = 9

$ws.Range("A169").Value = "Brunei"
$ws.Range("B169").Value = 141
$ws.Range("D169").Value = 138
$ws.Range("E169").Value = 1
$ws.Range("H169").Value = 2

# Bahamas
$ws.Range("D174").Value = 49
$ws.Range("E174").Value = 42

# Belice overtakes Santa Lucia -> rows 201/202 swap
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("D202").Value = 18
$ws.Range("H202").Value = 0

# Seychelles overtakes Montserrat -> rows 210/211 swap
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
